# Applies the weekly reshuffle of Femacal de La Calera - Coco price records:
# the Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) values are updated per
# row to reflect the new weekly ordering of the data (row 6 is untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44389
$ws.Range("M2").Value2 = 20
$ws.Range("N2").Value2 = 20000
$ws.Range("O2").Value2 = 20000
$ws.Range("P2").Value2 = 20000
$ws.Range("S2").Value2 = 1000

# Row 3
$ws.Range("D3").Value2 = 44294
$ws.Range("M3").Value2 = 25
$ws.Range("N3").Value2 = 25000
$ws.Range("O3").Value2 = 25000
$ws.Range("P3").Value2 = 25000
$ws.Range("S3").Value2 = 1250

# Row 4
$ws.Range("D4").Value2 = 44400
$ws.Range("M4").Value2 = 45
$ws.Range("N4").Value2 = 20000
$ws.Range("O4").Value2 = 20000
$ws.Range("P4").Value2 = 20000
$ws.Range("S4").Value2 = 1000

# Row 5
$ws.Range("D5").Value2 = 44307
$ws.Range("M5").Value2 = 30
$ws.Range("N5").Value2 = 22000
$ws.Range("O5").Value2 = 22000
$ws.Range("P5").Value2 = 22000
$ws.Range("S5").Value2 = 1100

# Row 7
$ws.Range("D7").Value2 = 44377
$ws.Range("M7").Value2 = 25
$ws.Range("N7").Value2 = 20000
$ws.Range("O7").Value2 = 20000
$ws.Range("P7").Value2 = 20000
$ws.Range("S7").Value2 = 1000

# Row 8
$ws.Range("D8").Value2 = 44292
$ws.Range("M8").Value2 = 30
$ws.Range("N8").Value2 = 25000
$ws.Range("O8").Value2 = 25000
$ws.Range("P8").Value2 = 25000
$ws.Range("S8").Value2 = 1250

# Row 9
$ws.Range("D9").Value2 = 44406
$ws.Range("M9").Value2 = 20
$ws.Range("N9").Value2 = 20000
$ws.Range("O9").Value2 = 20000
$ws.Range("P9").Value2 = 20000
$ws.Range("S9").Value2 = 1000

# Row 10
$ws.Range("D10").Value2 = 44305
$ws.Range("M10").Value2 = 20
$ws.Range("N10").Value2 = 22000
$ws.Range("O10").Value2 = 22000
$ws.Range("P10").Value2 = 22000
$ws.Range("S10").Value2 = 1100

# Row 11
$ws.Range("D11").Value2 = 44300
$ws.Range("M11").Value2 = 45
$ws.Range("N11").Value2 = 22000
$ws.Range("O11").Value2 = 22000
$ws.Range("P11").Value2 = 22000
$ws.Range("S11").Value2 = 1100

# Row 12
$ws.Range("D12").Value2 = 44403
$ws.Range("M12").Value2 = 50
$ws.Range("N12").Value2 = 20000
$ws.Range("O12").Value2 = 20000
$ws.Range("P12").Value2 = 20000
$ws.Range("S12").Value2 = 1000

# Row 13
$ws.Range("D13").Value2 = 44382
$ws.Range("M13").Value2 = 24
$ws.Range("N13").Value2 = 20000
$ws.Range("O13").Value2 = 20000
$ws.Range("P13").Value2 = 20000
$ws.Range("S13").Value2 = 1000

# Row 14
$ws.Range("D14").Value2 = 44298
$ws.Range("M14").Value2 = 65
$ws.Range("N14").Value2 = 22000
$ws.Range("O14").Value2 = 22000
$ws.Range("P14").Value2 = 22000
$ws.Range("S14").Value2 = 1100

# Row 15
$ws.Range("D15").Value2 = 44448
$ws.Range("M15").Value2 = 30
$ws.Range("N15").Value2 = 22000
$ws.Range("O15").Value2 = 22000
$ws.Range("P15").Value2 = 22000
$ws.Range("S15").Value2 = 1100

# Row 16
$ws.Range("D16").Value2 = 44376
$ws.Range("M16").Value2 = 38
$ws.Range("N16").Value2 = 20000
$ws.Range("O16").Value2 = 20000
$ws.Range("P16").Value2 = 20000
$ws.Range("S16").Value2 = 1000

# Row 17
$ws.Range("D17").Value2 = 44301
$ws.Range("M17").Value2 = 38
$ws.Range("N17").Value2 = 22000
$ws.Range("O17").Value2 = 22000
$ws.Range("P17").Value2 = 22000
$ws.Range("S17").Value2 = 1100

# Row 18
$ws.Range("D18").Value2 = 44291
$ws.Range("M18").Value2 = 70
$ws.Range("N18").Value2 = 25000
$ws.Range("O18").Value2 = 25000
$ws.Range("P18").Value2 = 25000
$ws.Range("S18").Value2 = 1250

# Row 19
$ws.Range("D19").Value2 = 44413
$ws.Range("M19").Value2 = 45
$ws.Range("N19").Value2 = 20000
$ws.Range("O19").Value2 = 20000
$ws.Range("P19").Value2 = 20000
$ws.Range("S19").Value2 = 1000

# Row 20
$ws.Range("D20").Value2 = 44385
$ws.Range("M20").Value2 = 36
$ws.Range("N20").Value2 = 20000
$ws.Range("O20").Value2 = 20000
$ws.Range("P20").Value2 = 20000
$ws.Range("S20").Value2 = 1000
